# Auto-generated edit script: updates cryptos Price (D) and Volume/1h (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.839.20"
$ws.Range("E2").Value = "  -0.13%  "

$ws.Range("D3").Value = "2.083.47"
$ws.Range("E3").Value = "  -0.40%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("E5").Value = "  -0.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.624"
$ws.Range("E6").Value = "  -0.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.05"
$ws.Range("E7").Value = "  +2.61%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  +1.52%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0786"
$ws.Range("E10").Value = "  +0.57%  "

$ws.Range("E11").Value = "  +1.41%  "

$ws.Range("D12").Value = "2.390.57"
$ws.Range("E12").Value = "  +0.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.73"
$ws.Range("E13").Value = "  +2.26%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.17"
$ws.Range("E14").Value = "  +0.10%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.774"
$ws.Range("E15").Value = "  +1.25%  "

$ws.Range("E16").Value = "  +1.98%  "

$ws.Range("D17").Value = "2.054.33"
$ws.Range("E17").Value = "  -1.56%  "

$ws.Range("D18").Value = "37.776.51"
$ws.Range("E18").Value = "  -0.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.13"
$ws.Range("E19").Value = "  -0.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.48"
$ws.Range("E20").Value = "  +0.88%  "

$ws.Range("D21").Value = "0.0₃0847"
$ws.Range("E21").Value = "  +3.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.90"
$ws.Range("E22").Value = "  -0.24%  "

$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("E24").Value = "  -0.60%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.40"
$ws.Range("E25").Value = "  +0.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.57"
$ws.Range("E26").Value = "  +0.42%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.23"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.138"
$ws.Range("E28").Value = "  -1.06%  "

$ws.Range("E29").Value = "  -1.81%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.47"
$ws.Range("E30").Value = "  -0.15%  "

$ws.Range("E31").Value = "  +1.42%  "

$ws.Range("E32").Value = "  +2.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.71"
$ws.Range("E33").Value = "  +2.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0631"
$ws.Range("E34").Value = "  +0.66%  "

$ws.Range("E35").Value = "  -0.51%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.43"
$ws.Range("E36").Value = "  +0.66%  "

$ws.Range("E37").Value = "  -0.98%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.40"
$ws.Range("E39").Value = "  -0.89%  "

$ws.Range("E40").Value = "  -2.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "99.27"
$ws.Range("E41").Value = "  +1.98%  "

$ws.Range("E42").Value = "  +2.08%  "

$ws.Range("E43").Value = "  -1.48%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.94"
$ws.Range("E44").Value = "  +7.50%  "

$ws.Range("D45").Value = "1.446.54"
$ws.Range("E45").Value = "  -0.45%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.15"
$ws.Range("E47").Value = "  +2.56%  "

$ws.Range("E48").Value = "  +0.37%  "

$ws.Range("E49").Value = "  -0.41%  "

$ws.Range("E50").Value = "  -0.39%  "

$ws.Range("D51").Value = "2.276.32"
$ws.Range("E51").Value = "  -0.19%  "
